$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 4 (pushes existing rows 4-22 down to 6-24)
$ws.Rows("4:5").Insert()

# Row 4: new weekly data (Membrillo Champion - Primera)
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45071
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104003
$ws.Range("J4").Value = "Membrillo"
$ws.Range("K4").Value = "Champion"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("Q4").Value = "$/caja 18 kilos empedrada"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 667
$ws.Range("T4").Value = 18

# Row 5: new weekly data (Membrillo Champion - Segunda)
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 45071
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104003
$ws.Range("J5").Value = "Membrillo"
$ws.Range("K5").Value = "Champion"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("Q5").Value = "$/caja 18 kilos empedrada"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 556
$ws.Range("T5").Value = 18
